# Fruta / hortaliza, semanal
# A new weekly price reading is added for "Terminal La Palmera de La Serena -
# Tuna": insert a row so the sheet grows from 8 to 9 data rows, then make
# row 5 the newest entry (updated Fecha + Volumen) while rows 6-9 hold the
# data that used to live in rows 5-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: shift old rows 6, 7, 8 down to 7, 8, 9.
$ws.Rows.Item(6).Insert()

# --- Row 5: newest week (date + volume changed, quality/prices unchanged) ---
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44603
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107011
$ws.Range("J5").Value = "Tuna"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 240
$ws.Range("N5").Value = 14500
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14750
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 819
$ws.Range("T5").Value = 18

# --- Row 6: previously row 5's data ---
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Terminal La Palmera de La Serena"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44294
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107011
$ws.Range("J6").Value = "Tuna"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 14500
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14750
$ws.Range("Q6").Value = "$/caja 18 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 819
$ws.Range("T6").Value = 18

# --- Row 7: previously row 6's data ---
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Terminal La Palmera de La Serena"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44294
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107011
$ws.Range("J7").Value = "Tuna"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 240
$ws.Range("N7").Value = 12500
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12750
$ws.Range("Q7").Value = "$/caja 18 kilos"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 708
$ws.Range("T7").Value = 18

# --- Row 8: previously row 7's data ---
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44294
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107011
$ws.Range("J8").Value = "Tuna"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 10500
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 10750
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 597
$ws.Range("T8").Value = 18

# --- Row 9 (new row): previously row 8's data ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44595
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107011
$ws.Range("J9").Value = "Tuna"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 15500
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15750
$ws.Range("Q9").Value = "$/caja 18 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 875
$ws.Range("T9").Value = 18
